# func ptr | practical example interop
#
# Slide 10 ("Function pointers vs delegates performance results" / "TBD")
# becomes the "Function pointers practical example" / "DEMO" slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- Title shape: "Function pointers vs delegates performance results" ---
# The title textbox is made up of two runs:
#   run 1 (chars 1-31):  "Function pointers vs delegates "
#   run 2 (chars 32-50): "performance results"
# Collapse them into a single run reading "Function pointers practical
# example" by clearing the second run's text and rewriting the first run.
$title = $s.Shapes.Item(4)
$titleRange = $title.TextFrame.TextRange
$run1 = $titleRange.Characters(1, 31)
$run2 = $titleRange.Characters(32, 19)
$run2.Text = ""
$run1.Text = "Function pointers practical example"

# Re-apply formatting across the whole (now single-run) range using a
# fresh TextRange reference so the edit doesn't re-split the run.
$titleRange = $title.TextFrame.TextRange
$titleRange.Font.Size = 24
$titleRange.Font.Bold = $true
$titleRange.Font.Color.RGB = 16777215

# The textbox has <a:spAutoFit/>, so collapsing the title to a single
# line re-shrinks its height automatically - no manual resize needed.

# --- Subtitle shape: "TBD" -> "DEMO" ---
$subtitle = $s.Shapes.Item(5)
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Text = "DEMO"
$subtitleRange.Font.Size = 24
$subtitleRange.Font.Bold = $true
$subtitleRange.Font.Color.RGB = 16777215
